$d = $word.ActiveDocument

# Helper: split the run at a given absolute character offset by adding then
# immediately deleting a temporary bookmark there. Word automatically breaks
# the run in two at that point and the split survives the bookmark removal.
function Split-RunAt($offset) {
    $rng = $d.Range($offset, $offset)
    $bmName = "TmpSplitMarker"
    if ($d.Bookmarks.Exists($bmName)) {
        $d.Bookmarks($bmName).Delete()
    }
    $d.Bookmarks.Add($bmName, $rng) | Out-Null
    $d.Bookmarks($bmName).Delete()
}

# ---------------------------------------------------------------------
# 1) Remove the original "_GoBack" bookmark (previously located right
#    before the "V" in "(Annexure-V: Contractors Claims ...)").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Turn "(Annexure-V: Contractors Claims and Their Settlement)" into
#    "(Annexure-V: Contractors Claims and Their Settlement, Annexure-VI:
#    Approved Time Extension by competent Authority)"
# ---------------------------------------------------------------------
$marker = "Contractors Claims and Their Settlement"
$text = $d.Content.Text
$idx = $text.IndexOf($marker)
$closeParenIdx = $idx + $marker.Length

$rngParen = $d.Range($closeParenIdx, $closeParenIdx + 1)
# sanity check - should currently contain ")"
if ($rngParen.Text -ne ")") {
    throw "Unexpected text at close paren location: [$($rngParen.Text)]"
}
$newTail = ", Annexure-VI: Approved Time Extension by competent Authority)"
$rngParen.Text = $newTail

# Now split that freshly inserted text (plus the pre-existing trailing space
# run that follows it, which gets swept up into the same merged run) into
# the individual runs seen in the target document: "," | " " | "Annexure-VI"
# | ":" | " Approved Time Extension by competent Authority)" | " "
$p1 = $closeParenIdx + 1   # after ","
$p2 = $p1 + 1              # after " "
$p3 = $p2 + "Annexure-VI".Length  # after "Annexure-VI"
$p4 = $p3 + 1              # after ":"
$p5 = $closeParenIdx + $newTail.Length  # after the final ")"

Split-RunAt $p5
Split-RunAt $p4
Split-RunAt $p3
Split-RunAt $p2
Split-RunAt $p1

# ---------------------------------------------------------------------
# 3) Split "may be" into "may " and "be", re-inserting the "_GoBack"
#    bookmark between the two new runs.
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idxMay = $text.IndexOf("may be")
$splitPoint = $idxMay + 4  # right between "may " and "be"

$rngBookmark = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $rngBookmark) | Out-Null
